$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.127.76'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.266.24'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '299.02'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.57'
$ws.Range('E6').Value = '  -3.90%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.496'
$ws.Range('E7').Value = '  -2.14%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('E10').Value = '  -3.41%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '47.75'
$ws.Range('E12').Value = '  -8.11%  '
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.67'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.618.44'
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.51'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.291.21'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.783'
$ws.Range('E18').Value = '  -4.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '42.057.95'
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.62'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0889'
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.98'
$ws.Range('E22').Value = '  -2.39%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '66.42'
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '234.86'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.96'
$ws.Range('E25').Value = '  -0.88%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('E27').Value = '  -3.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '23.97'
$ws.Range('E28').Value = '  -5.94%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.18'
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.57'
$ws.Range('E30').Value = '  +5.04%  '
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '33.51'
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('E34').Value = '  -3.18%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.47'
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('E36').Value = '  -5.03%  '
$ws.Range('E37').Value = '  -2.48%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0685'
$ws.Range('E38').Value = '  -4.52%  '
$ws.Range('E39').Value = '  -3.92%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0983'
$ws.Range('E40').Value = '  -2.04%  '
$ws.Range('E41').Value = '  -2.42%  '
$ws.Range('E42').Value = '  -5.07%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.43'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.961.90'
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('E45').Value = '  -1.39%  '
$ws.Range('E46').Value = '  -7.88%  '
$ws.Range('E47').Value = '  -6.44%  '
$ws.Range('E48').Value = '  -3.03%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.490.56'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '52.09'
$ws.Range('E50').Value = '  -5.71%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.49'
$ws.Range('E51').Value = '  -3.38%  '
